# Commit: "Test bench en C et deplacement de la documentation"
# The RISC-V compiler was re-run, producing updated machine-code encodings
# for several instructions in column A (and the "expected" column E), and
# two new "expected encoding" annotations appear (E3, E10) that did not
# exist before. Column B (mnemonics) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cells E3 / E10: copy formatting (style) from an existing E-column
# cell (style "2", same as column A) before writing their value, so no new
# style entries are introduced. ---
$ws.Range("A3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "111111100000 00010 000 00010 0010011"

$ws.Range("A10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "111111100100 00110 010 00010 0000011"

# --- Updated machine-code column A values (recompiled encodings) ---
$ws.Range("A4").Value  = "00000000011000010010111000100011"
$ws.Range("A7").Value  = "01111110001100110010001000100011"
$ws.Range("A9").Value  = "01111110001100110010011000100011"
$ws.Range("A10").Value = "11111110010000110010000100000011"
$ws.Range("A11").Value = "11111110110000110010000110000011"
$ws.Range("A15").Value = "01111110000000110010010000100011"
$ws.Range("A17").Value = "11111110110000110010000110000011"
$ws.Range("A19").Value = "01111110001100110010011000100011"
$ws.Range("A20").Value = "11111110100000110010000110000011"
$ws.Range("A22").Value = "01111110001100110010010000100011"
$ws.Range("A23").Value = "11111110100000110010000100000011"
$ws.Range("A28").Value = "11111110110000110010000110000011"
$ws.Range("A29").Value = "00000001110000010010001100000011"

# --- Updated "expected encoding" column E value ---
$ws.Range("E4").Value = "0000000 00110 00010 010 11100 0100011"

# --- Selection moved (user was last looking at F24) ---
$ws.Range("F24").Select()
